$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
# Row 32
$ws.Range("H32").Value = 1164.1578
$ws.Range("I32").Value = 452.375
$ws.Range("K32").Value = 452.375
$ws.Range("M32").Value = -126.375
# Row 62
$ws.Range("H62").Value = 1074.6522
$ws.Range("I62").Value = 1050.25
$ws.Range("J62").Value = 1101.2727
$ws.Range("K62").Value = 1050.25
$ws.Range("L62").Value = 1101.2727
$ws.Range("M62").Value = -426.25
$ws.Range("N62").Value = -2349.2727
# Row 65
$ws.Range("H65").Value = 1074.6522
$ws.Range("I65").Value = 1050.25
$ws.Range("J65").Value = 1101.2727
$ws.Range("K65").Value = 5251.25
$ws.Range("L65").Value = 5506.363499999999
$ws.Range("M65").Value = -2131.25
$ws.Range("N65").Value = -11746.3635
# Row 129
$ws.Range("H129").Value = 864.8099999999999
$ws.Range("I129").Value = 420
$ws.Range("J129").Value = 962.45123
$ws.Range("K129").Value = 1260
$ws.Range("L129").Value = 2887.35369
$ws.Range("M129").Value = 3740
$ws.Range("N129").Value = -12887.35369
# Row 132
$ws.Range("H132").Value = 961.13336
$ws.Range("I132").Value = 862.4583
$ws.Range("J132").Value = 1355.8334
$ws.Range("K132").Value = 2587.3749
$ws.Range("L132").Value = 4067.5002
$ws.Range("M132").Value = -57.3748999999998
$ws.Range("N132").Value = -9127.5002
# Row 137
$ws.Range("H137").Value = 1388.1034
$ws.Range("I137").Value = 1188.0476
$ws.Range("J137").Value = 1913.25
$ws.Range("K137").Value = 3564.142800000001
$ws.Range("L137").Value = 5739.75
$ws.Range("M137").Value = -1014.142800000001
$ws.Range("N137").Value = -10839.75
# Row 141
$ws.Range("H141").Value = 1562.8214
$ws.Range("I141").Value = 1509.5927
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 4528.7781
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 651.2219000000005
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 5500.875
$ws.Range("I26").Value = 4001.4
$ws.Range("K26").Value = 4001.4
$ws.Range("M26").Value = -3671.4
# Row 32
$ws.Range("H32").Value = 4361.537
$ws.Range("I32").Value = 3275.9473
$ws.Range("J32").Value = 10549.4
$ws.Range("K32").Value = 3275.9473
$ws.Range("L32").Value = 10549.4
$ws.Range("M32").Value = -2988.9473
$ws.Range("N32").Value = -11123.4
# Row 61
$ws.Range("H61").Value = 5604.88
$ws.Range("I61").Value = 6881.684
$ws.Range("J61").Value = 1561.6666
$ws.Range("K61").Value = 6881.684
$ws.Range("L61").Value = 1561.6666
$ws.Range("M61").Value = -6669.684
$ws.Range("N61").Value = -1985.6666
# Row 74
$ws.Range("H74").Value = 1181.8857
$ws.Range("I74").Value = 1248.6522
$ws.Range("J74").Value = 1053.9166
$ws.Range("K74").Value = 1248.6522
$ws.Range("L74").Value = 1053.9166
$ws.Range("M74").Value = -374.6522
$ws.Range("N74").Value = -2801.9166
# Row 77
$ws.Range("H77").Value = 1181.8857
$ws.Range("I77").Value = 1248.6522
$ws.Range("J77").Value = 1053.9166
$ws.Range("K77").Value = 6243.261
$ws.Range("L77").Value = 5269.583000000001
$ws.Range("M77").Value = -1875.261
$ws.Range("N77").Value = -14005.583
# Row 110
$ws.Range("H110").Value = 1285.7142
$ws.Range("I110").Value = 1366.6666
$ws.Range("J110").Value = 1225
$ws.Range("K110").Value = 1366.6666
$ws.Range("L110").Value = 1225
$ws.Range("M110").Value = 678.3334
$ws.Range("N110").Value = -5315
# Row 123
$ws.Range("H123").Value = 49424
$ws.Range("J123").Value = 49424
$ws.Range("L123").Value = 49424
$ws.Range("N123").Value = -59224
# Row 132
$ws.Range("H132").Value = 3867.7273
$ws.Range("I132").Value = 2267.8096
$ws.Range("J132").Value = 6667.5835
$ws.Range("K132").Value = 6803.4288
$ws.Range("L132").Value = 20002.7505
$ws.Range("M132").Value = -4273.4288
$ws.Range("N132").Value = -25062.7505
# Row 133
$ws.Range("H133").Value = 39000
$ws.Range("J133").Value = 39000
$ws.Range("L133").Value = 39000
$ws.Range("N133").Value = -44060
# Row 136
$ws.Range("H136").Value = 5604.88
$ws.Range("I136").Value = 6881.684
$ws.Range("J136").Value = 1561.6666
$ws.Range("K136").Value = 20645.052
$ws.Range("L136").Value = 4684.9998
$ws.Range("M136").Value = -18095.052
$ws.Range("N136").Value = -9784.9998

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 7247898.5
$ws.Range("I86").Value = 10418074
$ws.Range("J86").Value = 1783.8572
$ws.Range("K86").Value = 10418074
$ws.Range("L86").Value = 1783.8572
$ws.Range("M86").Value = -10416951
$ws.Range("N86").Value = -4029.8572
# Row 89
$ws.Range("H89").Value = 7247898.5
$ws.Range("I89").Value = 10418074
$ws.Range("J89").Value = 1783.8572
$ws.Range("K89").Value = 52090370
$ws.Range("L89").Value = 8919.286
$ws.Range("M89").Value = -52084754
$ws.Range("N89").Value = -20151.286
# Row 94
$ws.Range("H94").Value = 1063.7188
$ws.Range("I94").Value = 655.3461
$ws.Range("J94").Value = 2833.3333
$ws.Range("K94").Value = 655.3461
$ws.Range("L94").Value = 2833.3333
$ws.Range("M94").Value = -204.3461
$ws.Range("N94").Value = -3735.3333
# Row 99
$ws.Range("H99").Value = 200001570
$ws.Range("I99").Value = 333334300
$ws.Range("K99").Value = 333334300
$ws.Range("M99").Value = -333332802
# Row 107
$ws.Range("H107").Value = 1148.5454
$ws.Range("I107").Value = 1242.625
$ws.Range("J107").Value = 897.6667
$ws.Range("K107").Value = 1242.625
$ws.Range("L107").Value = 897.6667
$ws.Range("M107").Value = 677.375
$ws.Range("N107").Value = -4737.6667
# Row 134
$ws.Range("H134").Value = 4916.647
$ws.Range("I134").Value = 5413.5557
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 16240.6671
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -13705.6671
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3232.9756
$ws.Range("I31").Value = 1628.9474
$ws.Range("J31").Value = 4618.273
$ws.Range("K31").Value = 1628.9474
$ws.Range("L31").Value = 4618.273
$ws.Range("M31").Value = -1333.9474
$ws.Range("N31").Value = -5208.273
# Row 32
$ws.Range("H32").Value = 650
$ws.Range("I32").Value = 650
$ws.Range("K32").Value = 650
$ws.Range("M32").Value = -334
# Row 34
$ws.Range("H34").Value = 3232.9756
$ws.Range("I34").Value = 1628.9474
$ws.Range("J34").Value = 4618.273
$ws.Range("K34").Value = 1628.9474
$ws.Range("L34").Value = 4618.273
$ws.Range("M34").Value = -1426.9474
$ws.Range("N34").Value = -5022.273
# Row 58
$ws.Range("H58").Value = 1242.3256
$ws.Range("I58").Value = 909.9286
$ws.Range("J58").Value = 1862.8
$ws.Range("K58").Value = 909.9286
$ws.Range("L58").Value = 1862.8
$ws.Range("M58").Value = -706.9286
$ws.Range("N58").Value = -2268.8
# Row 132
$ws.Range("H132").Value = 2440.5173
$ws.Range("I132").Value = 2281
$ws.Range("J132").Value = 3052
$ws.Range("K132").Value = 6843
$ws.Range("L132").Value = 9156
$ws.Range("M132").Value = -4313
$ws.Range("N132").Value = -14216
# Row 134
$ws.Range("H134").Value = 1688.4736
$ws.Range("I134").Value = 1596.6428
$ws.Range("J134").Value = 1945.6
$ws.Range("K134").Value = 4789.928400000001
$ws.Range("L134").Value = 5836.799999999999
$ws.Range("M134").Value = -2254.928400000001
$ws.Range("N134").Value = -10906.8
# Row 136
$ws.Range("H136").Value = 1242.3256
$ws.Range("I136").Value = 909.9286
$ws.Range("J136").Value = 1862.8
$ws.Range("K136").Value = 2729.7858
$ws.Range("L136").Value = 5588.4
$ws.Range("M136").Value = -179.7857999999997
$ws.Range("N136").Value = -10688.4

$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 3300
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3300
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 9900
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -10466

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5812.516
$ws.Range("J70").Value = 5664.143
$ws.Range("L70").Value = 5664.143
$ws.Range("N70").Value = -6204.143
# Row 73
$ws.Range("H73").Value = 5812.516
$ws.Range("J73").Value = 5664.143
$ws.Range("L73").Value = 5664.143
$ws.Range("N73").Value = -7536.143
# Row 126
$ws.Range("H126").Value = 4912.3335
$ws.Range("I126").Value = 6985.222
$ws.Range("J126").Value = 1803
$ws.Range("K126").Value = 20955.666
$ws.Range("L126").Value = 5409
$ws.Range("M126").Value = -18485.666
$ws.Range("N126").Value = -10349
# Row 130
$ws.Range("H130").Value = 47983.332
$ws.Range("J130").Value = 47983.332
$ws.Range("L130").Value = 47983.332
$ws.Range("N130").Value = -58023.332
# Row 132
$ws.Range("H132").Value = 2919.4583
$ws.Range("I132").Value = 3576.8572
$ws.Range("J132").Value = 2408.1482
$ws.Range("K132").Value = 10730.5716
$ws.Range("L132").Value = 7224.444600000001
$ws.Range("M132").Value = -8200.571599999999
$ws.Range("N132").Value = -12284.4446

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 11637053
$ws.Range("I122").Value = 11909895
$ws.Range("J122").Value = 10000000
$ws.Range("K122").Value = 35729685
$ws.Range("L122").Value = 30000000
$ws.Range("M122").Value = -35727235
$ws.Range("N122").Value = -30004900
# Row 132
$ws.Range("H132").Value = 13100401
$ws.Range("I132").Value = 22923900
$ws.Range("J132").Value = 2401.2
$ws.Range("K132").Value = 68771700
$ws.Range("L132").Value = 7203.599999999999
$ws.Range("M132").Value = -68769170
$ws.Range("N132").Value = -12263.6
# Row 136
$ws.Range("H136").Value = 8239.321
$ws.Range("I136").Value = 10460.923
$ws.Range("J136").Value = 6313.933
$ws.Range("K136").Value = 31382.769
$ws.Range("L136").Value = 18941.799
$ws.Range("M136").Value = -28832.769
$ws.Range("N136").Value = -24041.799

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2311.4285
$ws.Range("I81").Value = 2108.7144
$ws.Range("J81").Value = 2514.1428
$ws.Range("K81").Value = 4217.4288
$ws.Range("L81").Value = 5028.2856
$ws.Range("M81").Value = -3156.4288
$ws.Range("N81").Value = -7150.2856
# Row 84
$ws.Range("H84").Value = 2311.4285
$ws.Range("I84").Value = 2108.7144
$ws.Range("J84").Value = 2514.1428
$ws.Range("K84").Value = 21087.144
$ws.Range("L84").Value = 25141.428
$ws.Range("M84").Value = -15783.144
$ws.Range("N84").Value = -35749.428
# Row 132
$ws.Range("H132").Value = 2102.8462
$ws.Range("I132").Value = 1426.6
$ws.Range("J132").Value = 2525.5
$ws.Range("K132").Value = 4279.799999999999
$ws.Range("L132").Value = 7576.5
$ws.Range("M132").Value = -1749.799999999999
$ws.Range("N132").Value = -12636.5
# Row 136
$ws.Range("H136").Value = 5200.923
$ws.Range("I136").Value = 6356.8887
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 19070.6661
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -16520.6661
$ws.Range("N136").Value = -12900
